# order cities by population (descending), keeping the header row (row 1) fixed
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:D13")
$sortKey   = $ws.Range("C2:C13")

# 2 = xlDescending
$dataRange.Sort($sortKey, 2)
